$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds date-like text (e.g. "2023-08-14"). Excel auto-converts such
# literals to date serials on plain .Value assignment, so the whole D range is
# pre-formatted as Text ("@") before any writes, then the formatting is cleared
# again afterwards so the new cells end up unstyled (same as the existing rows)
# while the values themselves remain literal text.
$ws.Range("D72:D97").NumberFormat = "@"

# Row 72
$ws.Range("A72").Value = 'Defense-gov Explore Feed'
$ws.Range("B72").Value = 'Military Commissions Media Invitation Announced for United States v. Abd al Hadi al Iraqi Pre-Sentencing Hearing'
$ws.Range("C72").Value = 'https://www.defense.gov/News/Advisories/Advisory/Article/3492079/military-commissions-media-invitation-announced-for-united-states-v-abd-al-hadi/'
$ws.Range("D72").Value = '2023-08-14'
$ws.Range("E72").Value = ' Department of Defense invites media to cover a pre-sentencing hearing in the case of United States v. Abd al-Hadi al-Iraqi, scheduled for Oct. 30 – Nov. 3.  The defendant entered a guilty plea in June of 2022.'

# Row 73
$ws.Range("A73").Value = 'Defense-gov Explore Feed'
$ws.Range("B73").Value = 'Secretary of Defense Lloyd J. Austin III Hosts Defense Minister Thabet al-Abbasi of Iraq'
$ws.Range("C73").Value = 'https://www.defense.gov/News/Advisories/Advisory/Article/3485025/secretary-of-defense-lloyd-j-austin-iii-hosts-defense-minister-thabet-al-abbasi/'
$ws.Range("D73").Value = '2023-08-07'
$ws.Range("E73").Value = 'retary of Defense Lloyd J. Austin III will host Defense Minister Thabet al-Abbasi of Iraq at an honor cordon ceremony at the Pentagon today.'

# Row 74
$ws.Range("A74").Value = 'Defense-gov Explore Feed'
$ws.Range("B74").Value = 'Secretary of Defense Lloyd J. Austin III Welcomes Mongolian Prime Minister Luvsannamsrain Oyun-Erdene to the Pentagon'
$ws.Range("C74").Value = 'https://www.defense.gov/News/Advisories/Advisory/Article/3480858/secretary-of-defense-lloyd-j-austin-iii-welcomes-mongolian-prime-minister-luvsa/'
$ws.Range("D74").Value = '2023-08-02'
$ws.Range("E74").Value = 'retary of Defense Lloyd J. Austin III will welcome Mongolian Prime Minister Luvsannamsrain Oyun-Erdene to the Pentagon in an enhanced honor cordon ceremony, August 3.'

# Row 75
$ws.Range("A75").Value = 'Defense-gov Explore Feed'
$ws.Range("B75").Value = 'Military Commissions Media Invitation Announced for United States v. Abd al-Rahim al-Nashiri Pre-Trial Hearing'
$ws.Range("C75").Value = 'https://www.defense.gov/News/Advisories/Advisory/Article/3480527/military-commissions-media-invitation-announced-for-united-states-v-abd-al-rahi/'
$ws.Range("D75").Value = '2023-08-02'
$ws.Range("E75").Value = ' Defense Department invites media to cover pre-trial proceedings for Abd al-Rahim al-Nashiri, the defendant charged in the attack on the USS Cole.'

# Row 76
$ws.Range("A76").Value = 'Defense-gov Explore Feed'
$ws.Range("B76").Value = 'Military Commissions Media Invitation Announced for United States v. Encep Nurjaman et al. Pre-Trial Hearing'
$ws.Range("C76").Value = 'https://www.defense.gov/News/Advisories/Advisory/Article/3480506/military-commissions-media-invitation-announced-for-united-states-v-encep-nurja/'
$ws.Range("D76").Value = '2023-08-02'
$ws.Range("E76").Value = ' Defense Department invites media to cover pre-trial proceedings for Encep Nurjaman, Mohammed Nazir Bin Lep and Mohammed Farik Bin Amin, all of whom have been charged jointly in connection with their alleged roles in the 2002 and 2003 bombings in Indonesia.'

# Row 77
$ws.Range("A77").Value = 'Defense-gov Explore Feed'
$ws.Range("B77").Value = 'On-Camera Press Briefing'
$ws.Range("C77").Value = 'https://www.defense.gov/News/Advisories/Advisory/Article/3478555/on-camera-press-briefing/'
$ws.Range("D77").Value = '2023-08-01'
$ws.Range("E77").Value = 'tagon Press Secretary Brig. Gen. Pat Ryder will hold an on-camera press briefing in the Pentagon Press Briefing Room.'

# Row 78
$ws.Range("A78").Value = 'Defense-gov Explore Feed'
$ws.Range("B78").Value = 'Deputy Secretary of Defense Kathleen Hicks Statement on the Release of the Commission on Planning, Programming, Budgeting, and Execution Reform Interim Report'
$ws.Range("C78").Value = 'https://www.defense.gov/News/Releases/Release/Article/3494248/deputy-secretary-of-defense-kathleen-hicks-statement-on-the-release-of-the-comm/'
$ws.Range("D78").Value = '2023-08-15'
$ws.Range("E78").Value = ' Defense Department must meet the urgency of today''s threats and tomorrow''s challenges with innovation in all portfolios.'

# Row 79
$ws.Range("A79").Value = 'Defense-gov Explore Feed'
$ws.Range("B79").Value = 'General Officer Assignments'
$ws.Range("C79").Value = 'https://www.defense.gov/News/Releases/Release/Article/3492361/general-officer-assignments/'
$ws.Range("D79").Value = '2023-08-14'
$ws.Range("E79").Value = ' acting chief of staff of the Army announced officer assignments.'

# Row 80
$ws.Range("A80").Value = 'Defense-gov Explore Feed'
$ws.Range("B80").Value = 'General Officer Announcement'
$ws.Range("C80").Value = 'https://www.defense.gov/News/Releases/Release/Article/3491979/general-officer-announcement/'
$ws.Range("D80").Value = '2023-08-14'
$ws.Range("E80").Value = 'retary of Defense Lloyd J. Austin III announced the president has made a nomination.'

# Row 81
$ws.Range("A81").Value = 'Defense-gov Explore Feed'
$ws.Range("B81").Value = 'Biden Administration Announces Additional Security Assistance for Ukraine'
$ws.Range("C81").Value = 'https://www.defense.gov/News/Releases/Release/Article/3491937/biden-administration-announces-additional-security-assistance-for-ukraine/'
$ws.Range("D81").Value = '2023-08-14'
$ws.Range("E81").Value = ' DOD announced additional security assistance to meet Ukraine''s critical security and defense needs. This announcement is the Biden Administration''s forty-fourth tranche of equipment to be provided for Ukraine since August 2021.'

# Row 82
$ws.Range("A82").Value = 'Defense-gov Explore Feed'
$ws.Range("B82").Value = 'Secretary of Defense Lloyd J. Austin III Statement on New Zealand''s Defense Policy and Strategy Statement and Future Force Design Principles'
$ws.Range("C82").Value = 'https://www.defense.gov/News/Releases/Release/Article/3490818/secretary-of-defense-lloyd-j-austin-iii-statement-on-new-zealands-defense-polic/'
$ws.Range("D82").Value = '2023-08-11'
$ws.Range("E82").Value = 'retary of Defense Lloyd J. Austin III issued a statement on New Zealand’s Defense Policy and Strategy Statement and Future Force Design Principles.'

# Row 83
$ws.Range("A83").Value = 'Defense-gov Explore Feed'
$ws.Range("B83").Value = 'DOD Announces Inaugural Innovation Challenge on Talent Management'
$ws.Range("C83").Value = 'https://www.defense.gov/News/Releases/Release/Article/3490776/dod-announces-inaugural-innovation-challenge-on-talent-management/'
$ws.Range("D83").Value = '2023-08-11'
$ws.Range("E83").Value = ' Office of the Undersecretary of Defense for Personnel and Readiness is hosting the inaugural Talent Management: From the Ground Up Innovation Challenge to capture the diversity of thought, experience, background and capability offered by our total force.'

# Row 84
$ws.Range("A84").Value = 'Defense-gov Explore Feed'
$ws.Range("B84").Value = 'Readout of Secretary of Defense Lloyd J. Austin III''s Phone Call With the President of Djibouti Ismail Omar Guelleh'
$ws.Range("C84").Value = 'https://www.defense.gov/News/Releases/Release/Article/3490095/readout-of-secretary-of-defense-lloyd-j-austin-iiis-phone-call-with-the-preside/'
$ws.Range("D84").Value = '2023-08-10'
$ws.Range("E84").Value = 'retary of Defense Lloyd J. Austin III spoke by phone with President of Djibouti Ismail Omar Guelleh to reaffirm the strength of the U.S.-Djibouti defense partnership.'

# Row 85
$ws.Range("A85").Value = 'Defense-gov Explore Feed'
$ws.Range("B85").Value = 'DOD Selects 2023-2024 Minerva-USIP Peace and Security Dissertation Fellows'
$ws.Range("C85").Value = 'https://www.defense.gov/News/Releases/Release/Article/3490020/dod-selects-2023-2024-minerva-usip-peace-and-security-dissertation-fellows/'
$ws.Range("D85").Value = '2023-08-10'
$ws.Range("E85").Value = ' Department of Defense today announced the 21 awardees of the 2023-2024 Minerva-U.S. Institute of Peace''s Peace and Security Dissertation Fellowship.'

# Row 86
$ws.Range("A86").Value = 'Defense-gov Explore Feed'
$ws.Range("B86").Value = 'DOD Announces Establishment of Generative AI Task Force'
$ws.Range("C86").Value = 'https://www.defense.gov/News/Releases/Release/Article/3489803/dod-announces-establishment-of-generative-ai-task-force/'
$ws.Range("D86").Value = '2023-08-10'
$ws.Range("E86").Value = ' DOD announced the establishment of a generative artificial intelligence task force to play a pivotal role in analyzing and integrating tools across the department.'

# Row 87
$ws.Range("A87").Value = 'Defense-gov Explore Feed'
$ws.Range("B87").Value = 'Readout of Acting Under Secretary of Defense for Policy Sasha Baker Phone Call With Saudi Deputy Minister of Defense, His Highness Prince Abdulrahman bin Ayyaf al-Muqrin'
$ws.Range("C87").Value = 'https://www.defense.gov/News/Releases/Release/Article/3488004/readout-of-acting-under-secretary-of-defense-for-policy-sasha-baker-phone-call/'
$ws.Range("D87").Value = '2023-08-09'
$ws.Range("E87").Value = 'ing Undersecretary of Defense for Policy Sasha Baker spoke with Saudi Arabia''s deputy minister of defense and reaffirmed the nation''s commitment to a strong defense partnership.'

# Row 88
$ws.Range("A88").Value = 'Defense-gov Explore Feed'
$ws.Range("B88").Value = 'Vets Visit San Francisco 49ers'' Training Camp'
$ws.Range("C88").Value = 'https://www.defense.gov/News/Feature-Stories/Story/Article/3492037/vets-visit-san-francisco-49ers-training-camp/'
$ws.Range("D88").Value = '2023-08-14'
$ws.Range("E88").Value = 'erans from Defense Logistics Agency Distribution San Joaquin, Calif., visited the San Francisco 49ers'' training camp through Operation: Care and Comfort, an organization supporting troops, veterans and their families across the globe.'

# Row 89
$ws.Range("A89").Value = 'Defense-gov Explore Feed'
$ws.Range("B89").Value = 'Medal of Honor Monday: Army Cpl. Lester Hammond Jr.'
$ws.Range("C89").Value = 'https://www.defense.gov/News/Feature-Stories/Story/Article/3488336/medal-of-honor-monday-army-cpl-lester-hammond-jr/'
$ws.Range("D89").Value = '2023-08-14'
$ws.Range("E89").Value = 'ing the Korean War, Army Cpl. Lester Hammond Jr. knew his patrol would be taken out by the enemy if he didn''t do something drastic. His decision saved soldiers'' lives, but not his own.'

# Row 90
$ws.Range("A90").Value = 'Defense-gov Explore Feed'
$ws.Range("B90").Value = 'Paratroopers Seek Jump Perfection at Leapfest'
$ws.Range("C90").Value = 'https://www.defense.gov/News/Feature-Stories/Story/Article/3487263/paratroopers-seek-jump-perfection-at-leapfest/'
$ws.Range("D90").Value = '2023-08-10'
$ws.Range("E90").Value = 'atroopers from the U.S. and around the globe joined together to test their skills and precision in the skies above Rhode Island during Leapfest 2023, the largest static line parachute training event and competition in the world.'

# Row 91
$ws.Range("A91").Value = 'Defense-gov Explore Feed'
$ws.Range("B91").Value = 'Young. Female. Scientist. Meet 4 of the Army''s Rising Civilian Stars'
$ws.Range("C91").Value = 'https://www.defense.gov/News/Feature-Stories/Story/Article/3486623/young-female-scientist-meet-4-of-the-armys-rising-civilian-stars/'
$ws.Range("D91").Value = '2023-08-08'
$ws.Range("E91").Value = 'the number of women joining the federal workforce in STEM grows, these four young scientists are making an impact in the Army''s Chemical, Biological, Radiological, Nuclear, Explosives Command.'

# Row 92
$ws.Range("A92").Value = 'Defense-gov Explore Feed'
$ws.Range("B92").Value = 'How These Scientists Are Like Google to Soldiers'
$ws.Range("C92").Value = 'https://www.defense.gov/News/Feature-Stories/Story/Article/3485570/how-these-scientists-are-like-google-to-soldiers/'
$ws.Range("D92").Value = '2023-08-08'
$ws.Range("E92").Value = 't the civilian Army scientists who protect soldiers by analyzing and identifying unknown samples to determine if they''re dangerous.'

# Row 93
$ws.Range("A93").Value = 'Defense-gov Explore Feed'
$ws.Range("B93").Value = 'Medal of Honor Monday: Army Capt. Loren D. Hagen'
$ws.Range("C93").Value = 'https://www.defense.gov/News/Feature-Stories/Story/Article/3480265/medal-of-honor-monday-army-capt-loren-d-hagen/'
$ws.Range("D93").Value = '2023-08-07'
$ws.Range("E93").Value = 'y Capt. Loren Douglas Hagen joined the Green Berets during the Vietnam War and even though he lost his life in battle,  the extraordinary heroism he displayed while leading his men during a harrowing mission earned him a posthumous Medal of Honor.'

# Row 94
$ws.Range("A94").Value = 'Defense-gov Explore Feed'
$ws.Range("B94").Value = 'Tomb Guard Braves Storm'
$ws.Range("C94").Value = 'https://www.defense.gov/News/Feature-Stories/Story/Article/3484065/tomb-guard-braves-storm/'
$ws.Range("D94").Value = '2023-08-04'
$ws.Range("E94").Value = 'y Pfc. Jessica Kwiatkowski kept her post in front of the Tomb of the Unknown Soldier at Arlington National Cemetery, Va., in the middle of a violent rainstorm.'

# Row 95
$ws.Range("A95").Value = 'Defense-gov Explore Feed'
$ws.Range("B95").Value = 'Marine Corps Parachute Rigger Stresses Safety, Readiness'
$ws.Range("C95").Value = 'https://www.defense.gov/News/Feature-Stories/Story/Article/3482116/marine-corps-parachute-rigger-stresses-safety-readiness/'
$ws.Range("D95").Value = '2023-08-03'
$ws.Range("E95").Value = 'ine Corps Sgt. Felix Lopez Saenz is a parachute rigger participating in Exercise Talisman Sabre 23 as a member of the 31st Marine Expeditionary Unit.'

# Row 96
$ws.Range("A96").Value = 'Defense-gov Explore Feed'
$ws.Range("B96").Value = 'Marine Corps Veteran Makes Australia Home'
$ws.Range("C96").Value = 'https://www.defense.gov/News/Feature-Stories/Story/Article/3481539/marine-corps-veteran-makes-australia-home/'
$ws.Range("D96").Value = '2023-08-03'
$ws.Range("E96").Value = 'n Seth Mooney served in the Marine Corps he vowed to one day return to Australia. Now, he runs a Tex-Mex barbecue business there.'

# Row 97
$ws.Range("A97").Value = 'Defense-gov Explore Feed'
$ws.Range("B97").Value = 'Rainbow of the Sea'
$ws.Range("C97").Value = 'https://www.defense.gov/News/Feature-Stories/Story/Article/3470367/rainbow-of-the-sea/'
$ws.Range("D97").Value = '2023-08-02'
$ws.Range("E97").Value = 'lors wear color-coded jerseys aboard Navy aircraft carriers to help fellow service member identify their roles as well as to ensure safety and order are maintained on the flight deck.'

$ws.Range("D72:D97").ClearFormats()
